# Update EUR->ARS rate: 2025-10-20T15:22:09Z
# Appends a new data row (date, time, rate) to the bottom of the sheet,
# matching the existing table's literal-text cell style (not Excel's
# auto-detected date/time values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$dateText = "2025-10-20"
$timeText = "15:22:09"
$rateText = "1.00 EUR = 1,703.8182"

# Column A holds an ISO-looking date string ("2025-10-20"). Typing that
# straight into a cell makes Excel's smart-entry treat it as a real date
# serial, which isn't what the source data represents (every existing row
# keeps it as plain text). Build it with a formula from an underscore
# variant (which Excel does NOT recognise as a date) so the computed
# result is plain text, then copy/paste-special just the *value* onto the
# target cell - this keeps the destination's default (unstyled) formatting
# intact, same as every other cell in the column.
$helperSource = $ws.Range("ZZ1")
$helperFormula = $ws.Range("ZZ2")

$helperSource.Value = $dateText.Replace("-", "_")
$helperFormula.Formula = "=SUBSTITUTE(ZZ1,""_"",""-"")"

$helperFormula.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)  # xlPasteValues

$ws.Range("ZZ1:ZZ2").Clear()
$excel.CutCopyMode = 0

# Columns B (time-of-day text) and C (rate text) don't look like dates, so
# a plain assignment already keeps them as literal text.
$ws.Range("B" + $newRow).Value = $timeText
$ws.Range("C" + $newRow).Value = $rateText
